$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new student's row (row 7)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "جاسم زعابي "
$ws.Range("C7").Value = 90
$ws.Range("D7").Value = 88
$ws.Range("E7").Value = 98

# Copy formatting from the row above so the new row matches existing style
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A7:E7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the row height used for the new row
$ws.Rows.Item(7).RowHeight = 20

# Update the active selection like Excel would after the edit
$ws.Range("H12").Select() | Out-Null
